$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SMFG")

$ws.Range("D8").Value = 19578800
$ws.Range("E8").Value = 17284700
$ws.Range("F8").Value = 16889600
$ws.Range("G8").Value = 17103100
$ws.Range("H8").Value = 16317400
$ws.Range("I8").Value = 15600500
$ws.Range("J8").Value = 15461400

$ws.Range("I15").Value = -1178700
$ws.Range("J15").Value = -1099400

$ws.Range("D17").Value = 7765300
$ws.Range("E17").Value = 6355200
$ws.Range("F17").Value = 4833400
$ws.Range("G17").Value = 3495000
$ws.Range("H17").Value = 2349600
$ws.Range("I17").Value = 2907000
$ws.Range("J17").Value = 2835200

$ws.Range("D18").Value = 11813400
$ws.Range("E18").Value = 10929500
$ws.Range("F18").Value = 12056100
$ws.Range("G18").Value = 13608100
$ws.Range("H18").Value = 13967800
$ws.Range("I18").Value = 12693500
$ws.Range("J18").Value = 12626200

$ws.Range("D20").Value = -1789400
$ws.Range("E20").Value = -2076600
$ws.Range("F20").Value = -3195400
$ws.Range("G20").Value = -1771300
$ws.Range("H20").Value = -1106600
$ws.Range("I20").Value = -4478700
$ws.Range("J20").Value = -4386900

$ws.Range("D21").Value = 12800200
$ws.Range("E21").Value = 11606500
$ws.Range("F21").Value = 11268200
$ws.Range("G21").Value = 14087600
$ws.Range("H21").Value = 14946800
$ws.Range("I21").Value = 10243400
$ws.Range("J21").Value = 9804500

$ws.Range("D23").Value = 10024000
$ws.Range("E23").Value = 8852900
$ws.Range("F23").Value = 8860700
$ws.Range("G23").Value = 11836800
$ws.Range("H23").Value = 12861200
$ws.Range("I23").Value = 8214800
$ws.Range("J23").Value = 8239300

$ws.Range("D24").Value = 2445500
$ws.Range("E24").Value = 1545400
$ws.Range("F24").Value = 2034400
$ws.Range("G24").Value = 3989900
$ws.Range("H24").Value = 4147600
$ws.Range("I24").Value = 2306600
$ws.Range("J24").Value = 4162400

$ws.Range("D26").Value = 7578500
$ws.Range("E26").Value = 7307500
$ws.Range("F26").Value = 6826300
$ws.Range("G26").Value = 7846900
$ws.Range("H26").Value = 8713600
$ws.Range("I26").Value = 5908200
$ws.Range("J26").Value = 4076900

$ws.Range("D27").Value = 6638700
$ws.Range("E27").Value = 6386900
$ws.Range("F27").Value = 5846100
$ws.Range("G27").Value = 6812600
$ws.Range("H27").Value = 7551600
$ws.Range("I27").Value = 3722700
$ws.Range("J27").Value = 3058400

$ws.Range("D32").Value = 1789400
$ws.Range("E32").Value = 2076600
$ws.Range("F32").Value = 3195400
$ws.Range("G32").Value = 1771300
$ws.Range("H32").Value = 1106600
$ws.Range("I32").Value = 4478700
$ws.Range("J32").Value = 4386900

$ws.Range("D33").Value = 6638700
$ws.Range("E33").Value = 6386900
$ws.Range("F33").Value = 5846100
$ws.Range("G33").Value = 6812600
$ws.Range("H33").Value = 7551600
$ws.Range("I33").Value = 3722700
$ws.Range("J33").Value = 3058400

$ws.Range("D35").Value = 6638700
$ws.Range("E35").Value = 6386900
$ws.Range("F35").Value = 5846100
$ws.Range("G35").Value = 6812600
$ws.Range("H35").Value = 7551600
$ws.Range("I35").Value = 3722700
$ws.Range("J35").Value = 3058400

$ws.Range("D41").Value = 485743000
$ws.Range("E41").Value = 423665000
$ws.Range("F41").Value = 386815000
$ws.Range("G41").Value = 359331000
$ws.Range("H41").Value = 298240000
$ws.Range("I41").Value = 216580000
$ws.Range("J41").Value = 84465200

$ws.Range("D42").Value = 212711000
$ws.Range("E42").Value = 220866000
$ws.Range("F42").Value = 214697000
$ws.Range("G42").Value = 200945000
$ws.Range("H42").Value = 161381000
$ws.Range("I42").Value = 182025000
$ws.Range("J42").Value = 39736300

$ws.Range("D47").Value = 6175700
$ws.Range("E47").Value = 5465100
$ws.Range("F47").Value = 5518200
$ws.Range("G47").Value = 5276800
$ws.Range("H47").Value = 3366300
$ws.Range("I47").Value = 4686000
$ws.Range("J47").Value = 1868200

$ws.Range("D48").Value = 31415200
$ws.Range("E48").Value = 28038800
$ws.Range("F48").Value = 26391600
$ws.Range("G48").Value = 25048500
$ws.Range("H48").Value = 21215000
$ws.Range("I48").Value = 36013800
$ws.Range("J48").Value = 20247600

$ws.Range("D49").Value = 7824900
$ws.Range("E49").Value = 8556400
$ws.Range("F49").Value = 7939500
$ws.Range("G49").Value = 7408800
$ws.Range("H49").Value = 7411900
$ws.Range("I49").Value = 15314900
$ws.Range("J49").Value = 16256900

$ws.Range("D52").Value = 3715700
$ws.Range("E52").Value = 3416400
$ws.Range("F52").Value = 2975100
$ws.Range("G52").Value = 4557000
$ws.Range("H52").Value = 2649700
$ws.Range("I52").Value = 7753800
$ws.Range("J52").Value = 8355100

$ws.Range("D54").Value = 1799400000
$ws.Range("E54").Value = 1788040000
$ws.Range("F54").Value = 1686740000
$ws.Range("G54").Value = 1658320000
$ws.Range("H54").Value = 1460270000
$ws.Range("I54").Value = 1335700000
$ws.Range("J54").Value = 1281880000

$ws.Range("D59").Value = 794700
$ws.Range("E59").Value = 727000
$ws.Range("F59").Value = 641100
$ws.Range("G59").Value = 693400
$ws.Range("H59").Value = 672000
$ws.Range("I59").Value = 4522000
$ws.Range("J59").Value = 2359800

$ws.Range("D61").Value = 179778000
$ws.Range("E61").Value = 171000000
$ws.Range("F61").Value = 141777000
$ws.Range("G61").Value = 145570000
$ws.Range("H61").Value = 109490000
$ws.Range("I61").Value = 58538900
$ws.Range("J61").Value = 94132200

$ws.Range("D62").Value = 6463100
$ws.Range("E62").Value = 5647200
$ws.Range("F62").Value = 6310700
$ws.Range("G62").Value = 7991400
$ws.Range("H62").Value = 3753600
$ws.Range("I62").Value = 5996900
$ws.Range("J62").Value = 1152100

$ws.Range("D66").Value = 1705450000
$ws.Range("E66").Value = 1700030000
$ws.Range("F66").Value = 1606130000
$ws.Range("G66").Value = 1576740000
$ws.Range("H66").Value = 1394450000
$ws.Range("I66").Value = 1277280000
$ws.Range("J66").Value = 1233310000

$ws.Range("D72").Value = 50220800
$ws.Range("E72").Value = 45563800
$ws.Range("F72").Value = 41017700
$ws.Range("G72").Value = 37070400
$ws.Range("H72").Value = 31476200
$ws.Range("I72").Value = 56970300
$ws.Range("J72").Value = 21765900

$ws.Range("D76").Value = 93955300
$ws.Range("E76").Value = 88004600
$ws.Range("F76").Value = 80606500
$ws.Range("G76").Value = 81581800
$ws.Range("H76").Value = 65820000
$ws.Range("I76").Value = 58417600
$ws.Range("J76").Value = 48567700

$ws.Range("D81").Value = 6638700
$ws.Range("E81").Value = 6386900
$ws.Range("F81").Value = 5846100
$ws.Range("G81").Value = 6812600
$ws.Range("H81").Value = 7551600
$ws.Range("I81").Value = 3722700
$ws.Range("J81").Value = 3058400

$ws.Range("D83").Value = 2773100
$ws.Range("E83").Value = 2750500
$ws.Range("F83").Value = 2404800
$ws.Range("G83").Value = 2248400
$ws.Range("H83").Value = 2083300
$ws.Range("I83").Value = 2026400
$ws.Range("J83").Value = 1563500

$ws.Range("D89").Value = 84473000
$ws.Range("E89").Value = 40836600
$ws.Range("F89").Value = -10190300
$ws.Range("G89").Value = 74015800
$ws.Range("H89").Value = 75065900
$ws.Range("I89").Value = 21455800
$ws.Range("J89").Value = 11941900

$ws.Range("D91").Value = -6441600
$ws.Range("E91").Value = -4482200
$ws.Range("F91").Value = -4784500
$ws.Range("G91").Value = -5233900
$ws.Range("H91").Value = -4204900
$ws.Range("I91").Value = -5244500
$ws.Range("J91").Value = -1175100

$ws.Range("D94").Value = -30693500
$ws.Range("E94").Value = 5255400
$ws.Range("F94").Value = 47378200
$ws.Range("G94").Value = -13173000
$ws.Range("H94").Value = 131265500
$ws.Range("I94").Value = 12727400
$ws.Range("J94").Value = -22347400

$ws.Range("D96").Value = -1975900
$ws.Range("E96").Value = -1853900
$ws.Range("F96").Value = -1916000
$ws.Range("G96").Value = -1545100
$ws.Range("H96").Value = -1536600
$ws.Range("I96").Value = -2444500
$ws.Range("J96").Value = -1283000

$ws.Range("D100").Value = -3168200
$ws.Range("E100").Value = -1505400
$ws.Range("F100").Value = -506200
$ws.Range("G100").Value = -2735400
$ws.Range("H100").Value = -9390900
$ws.Range("I100").Value = -6688300
$ws.Range("J100").Value = -2759000

$ws.Range("D101").Value = -848600
$ws.Range("E101").Value = -95400
$ws.Range("F101").Value = -900200
$ws.Range("G101").Value = 1606500
$ws.Range("H101").Value = 50500
$ws.Range("I101").Value = 5645000
$ws.Range("J101").Value = -557600

$ws.Range("D102").Value = 49762700
$ws.Range("E102").Value = 44491100
$ws.Range("F102").Value = 35781500
$ws.Range("G102").Value = 59713900
$ws.Range("H102").Value = 196991000
$ws.Range("I102").Value = 33139900
$ws.Range("J102").Value = -13722100
